$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: bold average of column J (|S*|/n) across instances
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary labels (col A) and bold formulas (col B)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting: B14:B17 bold, 12pt, vertically centered
$rngB = $ws.Range("B14:B17")
$rngB.Font.Bold = $true
$rngB.Font.Size = 12
$rngB.VerticalAlignment = -4108

# Formatting: J12 bold
$ws.Range("J12").Font.Bold = $true

# Row heights for the new summary rows
$ws.Range("A14:B17").RowHeight = 15.6

# Selection matching the author's final view state
$ws.Range("A14:B17").Select()
